$d = $word.ActiveDocument

# --- Change 1: split "Moved MarkUs, Piazza..." run into three runs with
#     proofErr spell-check markers bracketing "MarkUs" ---
$p1 = $d.Paragraphs(1)
$r1 = $p1.Range
$r1.MoveEnd(1, -1) # exclude the paragraph mark so the pPr is preserved
$xml1 = '<?xml version="1.0" encoding="UTF-8" standalone="yes"?><pkg:package xmlns:pkg="http://schemas.microsoft.com/office/2006/xmlPackage"><pkg:part pkg:name="/word/document.xml" pkg:contentType="application/vnd.openxmlformats-officedocument.wordprocessingml.document.main+xml"><pkg:xmlData><w:document xmlns:w="http://schemas.openxmlformats.org/wordprocessingml/2006/main"><w:body><w:p><w:r><w:t xml:space="preserve">Moved </w:t></w:r><w:proofErr w:type="spellStart"/><w:r><w:t>MarkUs</w:t></w:r><w:proofErr w:type="spellEnd"/><w:r><w:t>, Piazza into Resources to clear up navbar clutter</w:t></w:r></w:p></w:body></w:document></pkg:xmlData></pkg:part></pkg:package>'
$r1.InsertXML($xml1)

# --- Change 2: append the "What each team member worked on" section with
#     per-member percentage lines and bulleted task lines after the last
#     paragraph ---
$endPos = $d.Content.End
$rEnd = $d.Range($endPos, $endPos)
$xml2 = '<?xml version="1.0" encoding="UTF-8" standalone="yes"?><pkg:package xmlns:pkg="http://schemas.microsoft.com/office/2006/xmlPackage"><pkg:part pkg:name="/word/document.xml" pkg:contentType="application/vnd.openxmlformats-officedocument.wordprocessingml.document.main+xml"><pkg:xmlData><w:document xmlns:w="http://schemas.openxmlformats.org/wordprocessingml/2006/main"><w:body><w:p><w:r><w:t>What each team member worked on:</w:t></w:r></w:p><w:p><w:r><w:t>Arai: 33%</w:t></w:r></w:p><w:p><w:pPr><w:pStyle w:val="ListParagraph"/><w:numPr><w:ilvl w:val="0"/><w:numId w:val="1"/></w:numPr></w:pPr><w:r><w:t>Marks and remarks</w:t></w:r></w:p><w:p><w:r><w:t>Katarina: 33%</w:t></w:r></w:p><w:p><w:pPr><w:pStyle w:val="ListParagraph"/><w:numPr><w:ilvl w:val="0"/><w:numId w:val="1"/></w:numPr></w:pPr><w:r><w:t>Feedback</w:t></w:r></w:p><w:p><w:r><w:t>Dmitriy: 33%</w:t></w:r></w:p><w:p><w:pPr><w:pStyle w:val="ListParagraph"/><w:numPr><w:ilvl w:val="0"/><w:numId w:val="1"/></w:numPr></w:pPr><w:r><w:t>Sign up / log in system</w:t></w:r></w:p></w:body></w:document></pkg:xmlData></pkg:part></pkg:package>'
$rEnd.InsertXML($xml2)

Write-Host "Final paragraph count:" $d.Paragraphs.Count
for ($i = 1; $i -le $d.Paragraphs.Count; $i++) {
    Write-Host $i ":" $d.Paragraphs($i).Range.Text
}
